$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (MNN01 - Nguyen Quang Minh): DiaChi - drop "- Viet Nam" suffix
$ws.Range("G2").Value = "Làng sinh viên Hacinco - Nhân Chính - Thanh Xuân - Hà Nội"

# Row 5 (MNV03 - Pham Hong Nghia): Luong, NgayCMND, NoiCMND
$ws.Range("F5").Value = 1000000
$ws.Range("M5").Value = "12-12-2020"
$ws.Range("N5").Value = "Phú Thọ"

# Row 12 (MNV12 - Tran Hoang Anh): TenCV
$ws.Range("C12").Value = "Nhân viên"

# Row 13 (MNV24 - Phan Quoc Minh): DienThoai
$ws.Range("K13").Value = "0856865245"
